$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 22: reword cleaning-switch entry, add comment about backwards compatibility ---
$ws.Range("B22").Value = "Cleaning switch added"
$ws.Range("D22").Value = "This is logic controlled; so breaks backwards compatibility"

# --- Row 26: add new "Hardware" sub-heading (cell already carries the bold/left style) ---
$ws.Range("A26").Value = "Hardware"

# --- Rows 28-30: reword the hardware versioning rules ---
$ws.Range("A28").Value = "1st digit: Any change in hardware that breaks compatibility with firmware is a major revision and increments 1st digit"
$ws.Range("A29").Value = "2nd digit: Any change in sheet metal or manufactured components that affects anything else is a significant revision and increments 2nd digit"
$ws.Range("A30").Value = "3rd digit: Any change in purchased or manufactured components that affect nothing else is a minor revsion and increments 3rd digit"

# --- Row 31: blank separator row (keep the row present but empty) ---
$ws.Range("A31").Value = ""
$ws.Range("A31").Style = "Normal"

# --- Row 32: new bold "Software" sub-heading, mirroring the "Hardware" row's look ---
$ws.Range("B14").Copy($ws.Range("A32"))
$ws.Range("A32").Value = "Software"
$ws.Range("C25").Copy($ws.Range("C32"))

# --- Row 33: duplicate of the "Leading zero..." note, now under Software ---
$ws.Range("A33").Value = "Leading zero means pre-release. First release will be v1.0.0"
$ws.Range("A33").Style = "Normal"

# --- Row 34: now holds the software 1st-digit note; pick up the date-row look used in row 4 ---
$ws.Range("A34").Value = "1st digit same as first hardware version digit"
$ws.Range("A34").Style = "Normal"
$ws.Range("C4").Copy($ws.Range("C34"))
$ws.Range("C34").Value = ""

# --- Rows 35-36: new software versioning notes ---
$ws.Range("A35").Value = "2nd digit represents new features or major bugs"
$ws.Range("A35").Style = "Normal"
$ws.Range("A36").Value = "3rd digit is for minor hot fixes"
$ws.Range("A36").Style = "Normal"

# --- Row 38: blank row, just extends the used range/dimension ---
$ws.Range("A38").Value = ""
$ws.Range("A38").Style = "Normal"
